# Auto commit at 2025-12-29  8:34:20.71
# Append the next day's (2025-12-28, serial 46019) per-site readings for
# the two charging stations as two new rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
$row1 = $lastRow + 1
$row2 = $lastRow + 2

$ws.Cells.Item($row1, 1).Value = 46019
$ws.Cells.Item($row1, 2).Value = "四方坪站"
$ws.Cells.Item($row1, 3).Value = 9449.24
$ws.Cells.Item($row1, 4).Value = 8321.2900000000009
$ws.Cells.Item($row1, 5).Value = 3011.01
$ws.Cells.Item($row1, 6).Value = 386

$ws.Cells.Item($row2, 1).Value = 46019
$ws.Cells.Item($row2, 2).Value = "高岭站"
$ws.Cells.Item($row2, 3).Value = 5423.82
$ws.Cells.Item($row2, 4).Value = 4412.99
$ws.Cells.Item($row2, 5).Value = 1472.09
$ws.Cells.Item($row2, 6).Value = 182

# Match the author's saved viewport/selection state in the XLSX.
$ws.Range("H60").Select()
